# Update countries & provincias Spain
#
# The "Pais" sheet is sorted by total cases (column B) descending. This
# refresh re-sorts a handful of countries (Irak, Moldavia,
# Republica de Africa Central, Cabo Verde, Mozambique moved up a few
# spots each) and updates the day's case/death/recovery counters for
# the affected rows (plus a few unrelated rows whose totals simply grew).
# Since the row position/ranking is what encodes the sort order, each
# affected row gets both its country label (column A) and its stats
# (columns B-H) rewritten to the new post-refresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 2094237
$ws.Range("C4").Value = 4536
$ws.Range("D4").Value = 816467
$ws.Range("E4").Value = 1161638
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 98
$ws.Range("H4").Value = 116132

# Row 7
$ws.Range("B7").Value = 301579
$ws.Range("C7").Value = 3296
$ws.Range("D7").Value = 149767
$ws.Range("E7").Value = 143259
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 52
$ws.Range("H7").Value = 8553

# Row 20
$ws.Range("B20").Value = 97712
$ws.Range("C20").Value = 182
$ws.Range("D20").Value = 57960
$ws.Range("E20").Value = 31747
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 8005

# Row 31
$ws.Range("B31").Value = 41499
$ws.Range("C31").Value = 513
$ws.Range("D31").Value = 25946
$ws.Range("E31").Value = 15266
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 287

# Row 49
$ws.Range("A49").Value = "Irak"
$ws.Range("B49").Value = 17770
$ws.Range("C49").Value = 1095
$ws.Range("D49").Value = 6868
$ws.Range("E49").Value = 10406
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 39
$ws.Range("H49").Value = 496

# Row 50
$ws.Range("A50").Value = "Japon"
$ws.Range("B50").Value = 17292
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 15383
$ws.Range("E50").Value = 989
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 920

# Row 51
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 17269
$ws.Range("C51").Value = 602
$ws.Range("D51").Value = 11903
$ws.Range("E51").Value = 5330
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 36

# Row 52
$ws.Range("A52").Value = "Austria"
$ws.Range("B52").Value = 17064
$ws.Range("C52").Value = 30
$ws.Range("D52").Value = 15985
$ws.Range("E52").Value = 404
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 675

# Row 60
$ws.Range("A60").Value = "Moldavia"
$ws.Range("B60").Value = 11093
$ws.Range("C60").Value = 366
$ws.Range("D60").Value = 6229
$ws.Range("E60").Value = 4479
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 385

# Row 61
$ws.Range("A61").Value = "Ghana"
$ws.Range("B61").Value = 10856
$ws.Range("C61").Value = 498
$ws.Range("D61").Value = 3921
$ws.Range("E61").Value = 6887
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 48

# Row 98
$ws.Range("B98").Value = 2233
$ws.Range("C98").Value = 14
$ws.Range("D98").Value = 1902
$ws.Range("E98").Value = 247
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 84

# Row 100
$ws.Range("A100").Value = "Republica de Africa Central"
$ws.Range("B100").Value = 2044
$ws.Range("C100").Value = 92
$ws.Range("D100").Value = 360
$ws.Range("E100").Value = 1677
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 7

# Row 101
$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 1976
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 1153
$ws.Range("E101").Value = 815
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 8

# Row 102
$ws.Range("A102").Value = "Estonia"
$ws.Range("B102").Value = 1970
$ws.Range("C102").Value = 5
$ws.Range("D102").Value = 1703
$ws.Range("E102").Value = 198
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 69

# Row 103
$ws.Range("B103").Value = 1878
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 1196
$ws.Range("E103").Value = 671
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 11

# Row 124
$ws.Range("B124").Value = 1093
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 995
$ws.Range("E124").Value = 49
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 49

# Row 130
$ws.Range("B130").Value = 853
$ws.Range("C130").Value = 1
$ws.Range("D130").Value = 781
$ws.Range("E130").Value = 21
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 51

# Row 136
$ws.Range("A136").Value = "Cabo Verde"
$ws.Range("B136").Value = 697
$ws.Range("C136").Value = 40
$ws.Range("D136").Value = 294
$ws.Range("E136").Value = 397
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 6

# Row 137
$ws.Range("A137").Value = "San Marino"
$ws.Range("B137").Value = 694
$ws.Range("C137").Value = 3
$ws.Range("D137").Value = 520
$ws.Range("E137").Value = 132
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 42

# Row 138
$ws.Range("A138").Value = "Uganda"
$ws.Range("B138").Value = 686
$ws.Range("C138").Value = 7
$ws.Range("D138").Value = 161
$ws.Range("E138").Value = 525
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# Row 144
$ws.Range("A144").Value = "Mozambique"
$ws.Range("B144").Value = 509
$ws.Range("C144").Value = 20
$ws.Range("D144").Value = 145
$ws.Range("E144").Value = 362
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 2

# Row 145
$ws.Range("A145").Value = "Tanzania"
$ws.Range("B145").Value = 509
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 183
$ws.Range("E145").Value = 305
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 21

# Row 146
$ws.Range("A146").Value = "Ruanda"
$ws.Range("B146").Value = 494
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 313
$ws.Range("E146").Value = 179
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 2

# Row 147
$ws.Range("B147").Value = 489
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 414
$ws.Range("E147").Value = 72
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 3

# Row 148
$ws.Range("B148").Value = 488
$ws.Range("C148").Value = 1
$ws.Range("D148").Value = 460
$ws.Range("E148").Value = 27
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1
